$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 135, shifting existing rows 135..231 down to 136..232
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new record's data
$ws.Cells.Item(135, 1).Value = 6
$ws.Cells.Item(135, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(135, 3).Value = "Metropolitana"
$ws.Cells.Item(135, 4).Value = 44767
$ws.Cells.Item(135, 5).Value = 13
$ws.Cells.Item(135, 6).Value = 100112022
$ws.Cells.Item(135, 7).Value = "Arveja Verde"
$ws.Cells.Item(135, 8).Value = "Perfection"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 200
$ws.Cells.Item(135, 11).Value = 37000
$ws.Cells.Item(135, 12).Value = 40000
$ws.Cells.Item(135, 13).Value = 38200
$ws.Cells.Item(135, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(135, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(135, 16).Value = 1528
$ws.Cells.Item(135, 17).Value = 25
$ws.Cells.Item(135, 18).Value = "Hortaliza"
